$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table currently has 4 columns (Sum Sq, Df, F value, Pr(>F)) and
# 5 rows (1 header + 4 data rows). The edit switches the summary
# statistic from an F-test to a Chi-square test: the "F value" column
# is dropped, the last (all-NA) data row is dropped, and several
# labels/values are updated to match the new Anova (Chisq) output.

# Drop the "F value" column (3rd column).
$t.Columns.Item(3).Delete()

# Drop the last row (0.035 / 24 / NA / NA).
$t.Rows.Item($t.Rows.Count).Delete()

# Update the header row.
$t.Cell(1, 1).Range.Text = "Chisq"
$t.Cell(1, 3).Range.Text = "Pr(>Chisq)"

# Update the data rows (column 2 "Df" values are unchanged).
$t.Cell(2, 1).Range.Text = "8.874"
$t.Cell(2, 3).Range.Text = "0.012"

$t.Cell(3, 1).Range.Text = "1.966"
$t.Cell(3, 3).Range.Text = "0.161"

$t.Cell(4, 1).Range.Text = "0.463"
$t.Cell(4, 3).Range.Text = "0.793"
